$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.832.06"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "2.678.14"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D5").Value = "597.94"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "175.03"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "2.680.52"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "4.97"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").Value = "3.170.57"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("D16").Value = "71.806.34"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "26.14"
$ws.Range("E17").Value = "  -3.59%  "
$ws.Range("D18").Value = "2.673.60"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "370.96"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").Value = "2.02"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").Value = "71.99"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "4.32"
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("D27").Value = "9.74"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "0.0₃0968"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").Value = "502.53"
$ws.Range("E32").Value = "  -8.52%  "
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "163.50"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "19.51"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("E41").Value = "  -5.37%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").Value = "156.33"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "39.49"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "0.560"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "3.71"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  -0.55%  "
